$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-25 05:18:18'
$ws.Range('H2').Value = "'60%"
$ws.Range('N2').Value = '0.6 °C 4:58 TU'
$ws.Range('O2').Value = '1.6 °C'
$ws.Range('E3').Value = '2026-02-25 05:18:20'
$ws.Range('E4').Value = '2026-02-25 05:18:22'
$ws.Range('J4').Value = '1019.6 hPa'
$ws.Range('N4').Value = '1.5 °C 4:59 TU'
$ws.Range('O4').Value = '2.7 °C'
$ws.Range('E5').Value = '2026-02-25 05:18:25'
$ws.Range('E6').Value = '2026-02-25 05:18:27'
$ws.Range('J6').Value = '1019.1 hPa'
$ws.Range('N6').Value = '6.1 °C 4:57 TU'
$ws.Range('O6').Value = '8.5 °C'
$ws.Range('E7').Value = '2026-02-25 05:18:29'
$ws.Range('H7').Value = "'95%"
$ws.Range('J7').Value = '1018.7 hPa'
$ws.Range('K7').Value = '-0.1 MJ/m2'
$ws.Range('L7').Value = '11.2 km/h - 68º 4:58 TU'
$ws.Range('N7').Value = '10.0 °C 4:37 TU'
$ws.Range('O7').Value = '10.9 °C'
$ws.Range('E8').Value = '2026-02-25 05:18:32'
$ws.Range('H8').Value = "'48%"
$ws.Range('J8').Value = '1018.0 hPa'
$ws.Range('N8').Value = '12.7 °C 4:59 TU'
$ws.Range('O8').Value = '15.0 °C'
$ws.Range('E9').Value = '2026-02-25 05:18:34'
$ws.Range('N9').Value = '3.3 °C 4:59 TU'
$ws.Range('O9').Value = '5.2 °C'
$ws.Range('E10').Value = '2026-02-25 05:18:36'
$ws.Range('H10').Value = "'99%"
$ws.Range('E11').Value = '2026-02-25 05:18:39'
$ws.Range('H11').Value = "'88%"
$ws.Range('N11').Value = '1.8 °C 4:59 TU'
$ws.Range('O11').Value = '3.0 °C'
$ws.Range('E12').Value = '2026-02-25 05:18:41'
$ws.Range('O12').Value = '5.3 °C'
$ws.Range('E13').Value = '2026-02-25 05:18:43'
$ws.Range('J13').Value = '1026.6 hPa'
$ws.Range('K13').Value = '-0.1 MJ/m2'
$ws.Range('N13').Value = '-2.8 °C 4:56 TU'
$ws.Range('O13').Value = '-1.5 °C'
$ws.Range('E14').Value = '2026-02-25 05:18:45'
$ws.Range('N14').Value = '3.2 °C 4:36 TU'
$ws.Range('O14').Value = '5.5 °C'
$ws.Range('E15').Value = '2026-02-25 05:18:48'
$ws.Range('N15').Value = '4.1 °C 4:59 TU'
$ws.Range('O15').Value = '5.5 °C'
$ws.Range('E16').Value = '2026-02-25 05:18:50'
$ws.Range('E17').Value = '2026-02-25 05:18:52'
$ws.Range('H17').Value = "'25%"
$ws.Range('L17').Value = '33.1 km/h - 266º 4:52 TU'
$ws.Range('N17').Value = '7.5 °C 4:14 TU'
$ws.Range('O17').Value = '9.0 °C'
$ws.Range('E18').Value = '2026-02-25 05:18:54'
$ws.Range('J18').Value = '1019.5 hPa'
$ws.Range('N18').Value = '4.9 °C 4:38 TU'
$ws.Range('O18').Value = '6.2 °C'
$ws.Range('E19').Value = '2026-02-25 05:18:57'
$ws.Range('K19').Value = '-0.1 MJ/m2'
$ws.Range('L19').Value = '7.6 km/h - 81º 4:58 TU'
$ws.Range('O19').Value = '9.6 °C'
$ws.Range('E20').Value = '2026-02-25 05:18:59'
$ws.Range('H20').Value = "'43%"
$ws.Range('K20').Value = '-0.1 MJ/m2'
$ws.Range('O20').Value = '2.8 °C'
$ws.Range('E21').Value = '2026-02-25 05:19:01'
$ws.Range('H21').Value = "'75%"
$ws.Range('J21').Value = '1023.5 hPa'
$ws.Range('K21').Value = '-0.1 MJ/m2'
$ws.Range('N21').Value = '1.8 °C 4:49 TU'
$ws.Range('O21').Value = '3.6 °C'
$ws.Range('E22').Value = '2026-02-25 05:19:04'
$ws.Range('K22').Value = '-0.1 MJ/m2'
$ws.Range('E23').Value = '2026-02-25 05:19:06'
$ws.Range('H23').Value = "'26%"
$ws.Range('E24').Value = '2026-02-25 05:19:08'
$ws.Range('J24').Value = '1019.8 hPa'
$ws.Range('E25').Value = '2026-02-25 05:19:11'
$ws.Range('E26').Value = '2026-02-25 05:19:13'
$ws.Range('J26').Value = '1019.0 hPa'
$ws.Range('O26').Value = '9.2 °C'
$ws.Range('E27').Value = '2026-02-25 05:19:15'
$ws.Range('K27').Value = '-0.1 MJ/m2'
$ws.Range('E28').Value = '2026-02-25 05:19:18'
$ws.Range('J28').Value = '1020.5 hPa'
$ws.Range('O28').Value = '3.8 °C'
$ws.Range('E29').Value = '2026-02-25 05:19:20'
$ws.Range('L29').Value = '9.4 km/h - 351º 4:52 TU'
$ws.Range('M29').Value = '10.9 °C 4:56 TU'
$ws.Range('O29').Value = '9.1 °C'
$ws.Range('E30').Value = '2026-02-25 05:19:22'
$ws.Range('J30').Value = '1019.4 hPa'
$ws.Range('N30').Value = '6.5 °C 4:59 TU'
$ws.Range('O30').Value = '7.6 °C'
$ws.Range('E31').Value = '2026-02-25 05:19:25'
$ws.Range('H31').Value = "'92%"
$ws.Range('J31').Value = '1018.6 hPa'
$ws.Range('K31').Value = '-0.1 MJ/m2'
$ws.Range('E32').Value = '2026-02-25 05:19:27'
$ws.Range('H32').Value = "'73%"
$ws.Range('K32').Value = '-0.1 MJ/m2'
$ws.Range('N32').Value = '0.7 °C 4:54 TU'
$ws.Range('O32').Value = '2.3 °C'
$ws.Range('E33').Value = '2026-02-25 05:19:29'
$ws.Range('J33').Value = '1024.0 hPa'
$ws.Range('O33').Value = '2.3 °C'
$ws.Range('E34').Value = '2026-02-25 05:19:32'
$ws.Range('N34').Value = '0.1 °C 4:57 TU'
$ws.Range('E35').Value = '2026-02-25 05:19:34'
$ws.Range('H35').Value = "'35%"
$ws.Range('K35').Value = '-0.1 MJ/m2'
$ws.Range('O35').Value = '9.9 °C'
$ws.Range('E36').Value = '2026-02-25 05:19:36'
$ws.Range('J36').Value = '1019.1 hPa'
$ws.Range('M36').Value = '12.0 °C 4:34 TU'
$ws.Range('O36').Value = '8.7 °C'
$ws.Range('E37').Value = '2026-02-25 05:19:39'
$ws.Range('N37').Value = '-0.4 °C 4:44 TU'
$ws.Range('O37').Value = '1.0 °C'
$ws.Range('E38').Value = '2026-02-25 05:19:41'
$ws.Range('K38').Value = '-0.1 MJ/m2'
$ws.Range('E39').Value = '2026-02-25 05:19:43'
$ws.Range('N39').Value = '0.2 °C 4:31 TU'
$ws.Range('O39').Value = '2.0 °C'
$ws.Range('E40').Value = '2026-02-25 05:19:45'
$ws.Range('J40').Value = '1024.4 hPa'
$ws.Range('O40').Value = '1.5 °C'
$ws.Range('E41').Value = '2026-02-25 05:19:48'
$ws.Range('L41').Value = '6.5 km/h - 109º 4:58 TU'
$ws.Range('O41').Value = '8.9 °C'
$ws.Range('E42').Value = '2026-02-25 05:19:50'
$ws.Range('I42').Value = '0.1 mm'
$ws.Range('O42').Value = '8.2 °C'
$ws.Range('E43').Value = '2026-02-25 05:19:52'
$ws.Range('O43').Value = '3.7 °C'
$ws.Range('E44').Value = '2026-02-25 05:19:54'
$ws.Range('H44').Value = "'43%"
$ws.Range('K44').Value = '-0.1 MJ/m2'
$ws.Range('E45').Value = '2026-02-25 05:19:56'
$ws.Range('K45').Value = '-0.1 MJ/m2'
$ws.Range('N45').Value = '5.1 °C 4:33 TU'
$ws.Range('O45').Value = '6.5 °C'
$ws.Range('E46').Value = '2026-02-25 05:19:59'
$ws.Range('K46').Value = '-0.1 MJ/m2'
$ws.Range('N46').Value = '2.1 °C 4:44 TU'
$ws.Range('O46').Value = '3.6 °C'
